# Generate Report for Archive
# - Update status text "Ready for handoff" -> "In Translation" wherever it
#   appears (Overview!E2:F2/E3:F3 and the Status column on each locale sheet).
# - Shrink the columns that were auto-fit to the old (longer) status text
#   down to the new, narrower width.

$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    $used = $ws.UsedRange
    foreach ($cell in $used.Cells) {
        $val = [string]$cell.Value2
        if ($val -eq "Ready for handoff") {
            $cell.Value2 = "In Translation"
        }
    }
}

# NOTE: ColumnWidth is stored (and round-tripped through the xlsx) in whole
# pixels, so it always lands on a multiple of 1/6 here; 12.5 is the input
# that lands closest to the target raw column width (~13.41 character
# units) from the original diff.
$overview = $wb.Worksheets.Item("Overview")
$overview.Columns.Item(5).ColumnWidth = 12.5
$overview.Columns.Item(6).ColumnWidth = 12.5

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Columns.Item(3).ColumnWidth = 12.5

$dede = $wb.Worksheets.Item("de-de")
$dede.Columns.Item(3).ColumnWidth = 12.5
